$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E91").Select()
$excel.ActiveWindow.ScrollRow = 56
